$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Crime-complaint statistics table updates (rows 15-30) ---
    # Row 15
    $ws.Range("D15").NumberFormat = '#,##0'
    $ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G15").NumberFormat = '#,##0'
    $ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("D15").Value = 1
    $ws.Range("E15").Value = -100
    $ws.Range("G15").Value = 1
    $ws.Range("H15").Value = 0
    $ws.Range("J15").Value = 9
    $ws.Range("K15").Value = 0
    $ws.Range("N15").Value = -59.090909090909
    # Row 16
    $ws.Range("C16").NumberFormat = '#,##0'
    $ws.Range("C16").Value = 3
    $ws.Range("D16").Value = 3
    $ws.Range("E16").Value = 0
    $ws.Range("G16").Value = 14
    $ws.Range("H16").Value = -42.857142857142
    $ws.Range("I16").Value = 56
    $ws.Range("J16").Value = 81
    $ws.Range("K16").Value = -30.864197530864
    $ws.Range("L16").Value = -51.724137931034
    $ws.Range("M16").Value = -21.126760563380
    $ws.Range("N16").Value = -85.066666666666
    # Row 17
    $ws.Range("C17").Value = 3
    $ws.Range("D17").Value = 4
    $ws.Range("E17").Value = -25
    $ws.Range("F17").Value = 16
    $ws.Range("G17").Value = 20
    $ws.Range("H17").Value = -20
    $ws.Range("I17").Value = 82
    $ws.Range("J17").Value = 113
    $ws.Range("K17").Value = -27.433628318584
    $ws.Range("L17").Value = -6.818181818181
    $ws.Range("M17").Value = 17.142857142857
    $ws.Range("N17").Value = -68.093385214007
    # Row 18
    $ws.Range("C18").Value = 3
    $ws.Range("D18").Value = 5
    $ws.Range("E18").Value = -40
    $ws.Range("F18").Value = 18
    $ws.Range("G18").Value = 21
    $ws.Range("H18").Value = -14.285714285714
    $ws.Range("I18").Value = 80
    $ws.Range("J18").Value = 129
    $ws.Range("K18").Value = -37.984496124031
    $ws.Range("L18").Value = -55.801104972375
    $ws.Range("M18").Value = -23.076923076923
    $ws.Range("N18").Value = -77.715877437325
    # Row 19
    $ws.Range("D19").Value = 24
    $ws.Range("E19").Value = -54.166666666666
    $ws.Range("F19").Value = 59
    $ws.Range("G19").Value = 81
    $ws.Range("H19").Value = -27.160493827160
    $ws.Range("I19").Value = 334
    $ws.Range("J19").Value = 460
    $ws.Range("K19").Value = -27.391304347826
    $ws.Range("L19").Value = -23.569794050343
    $ws.Range("M19").Value = -2.052785923753
    $ws.Range("N19").Value = -50.445103857566
    # Row 20
    $ws.Range("C20").Value = 2
    $ws.Range("D20").Value = 3
    $ws.Range("E20").Value = -33.333333333333
    $ws.Range("F20").Value = 5
    $ws.Range("G20").Value = 4
    $ws.Range("H20").Value = 25
    $ws.Range("I20").Value = 17
    $ws.Range("J20").Value = 18
    $ws.Range("K20").Value = -5.555555555555
    $ws.Range("L20").Value = -22.727272727272
    $ws.Range("M20").Value = -19.047619047619
    $ws.Range("N20").Value = -92.672413793103
    # Row 21
    $ws.Range("C21").Value = 22
    $ws.Range("D21").Value = 40
    $ws.Range("E21").Value = -45
    $ws.Range("F21").Value = 107
    $ws.Range("G21").Value = 142
    $ws.Range("H21").Value = -24.647887323943
    $ws.Range("I21").Value = 578
    $ws.Range("J21").Value = 811
    $ws.Range("K21").Value = -28.729963008631
    $ws.Range("L21").Value = -32.712456344586
    $ws.Range("M21").Value = -5.709624796084
    $ws.Range("N21").Value = -69.95841995842
    # Row 22
    $ws.Range("G22").Value = 1
    $ws.Range("H22").Value = 0
    $ws.Range("L22").Value = -37.5
    # Row 23
    $ws.Range("G23").Value = 10
    $ws.Range("H23").Value = -10
    $ws.Range("I23").Value = 34
    $ws.Range("J23").Value = 63
    $ws.Range("K23").Value = -46.031746031746
    $ws.Range("L23").Value = -49.253731343283
    $ws.Range("M23").Value = -34.615384615384
    # Row 24
    $ws.Range("C24").Value = 36
    $ws.Range("D24").Value = 31
    $ws.Range("E24").Value = 16.129032258064
    $ws.Range("F24").Value = 136
    $ws.Range("G24").Value = 100
    $ws.Range("H24").Value = 36
    $ws.Range("I24").Value = 675
    $ws.Range("J24").Value = 653
    $ws.Range("K24").Value = 3.369065849923
    $ws.Range("L24").Value = -34.338521400778
    $ws.Range("M24").Value = -8.163265306122
    # Row 25
    $ws.Range("D25").Value = 14
    $ws.Range("E25").Value = 35.714285714285
    $ws.Range("F25").Value = 79
    $ws.Range("G25").Value = 52
    $ws.Range("H25").Value = 51.923076923076
    $ws.Range("I25").Value = 389
    $ws.Range("J25").Value = 380
    $ws.Range("K25").Value = 2.368421052631
    $ws.Range("L25").Value = -46.047156726768
    # Row 26
    $ws.Range("C26").Value = 5
    $ws.Range("D26").Value = 10
    $ws.Range("E26").Value = -50
    $ws.Range("F26").Value = 36
    $ws.Range("H26").Value = -2.702702702702
    $ws.Range("I26").Value = 176
    $ws.Range("J26").Value = 219
    $ws.Range("K26").Value = -19.634703196347
    $ws.Range("L26").Value = -16.981132075471
    $ws.Range("M26").Value = -16.190476190476
    # Row 27
    $ws.Range("D27").NumberFormat = '#,##0'
    $ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("D27").Value = 1
    $ws.Range("E27").Value = -100
    $ws.Range("G27").Value = 2
    $ws.Range("H27").Value = 0
    $ws.Range("J27").Value = 17
    $ws.Range("K27").Value = -35.294117647058
    # Row 28
    $ws.Range("F28").Value = 3
    $ws.Range("G28").Value = 2
    $ws.Range("H28").Value = 50
    $ws.Range("L28").Value = -32.558139534883
    # Row 29
    $ws.Range("G29").Value = 1
    # Row 30
    $ws.Range("G30").Value = 1
